# Update the environmental-data summary table after recalculating the
# Shannon diversity analysis (see commit message: "did the analysis
# again after calculating shannon diversity").
#
# Table 1 layout (1-based Word Table/Cell indices):
#   Row 10 "16S rRNA (Shannon Diversity)" -> col 13 "Landuse" p-value:
#                                             0.016 -> 0.023
#   Row 11 "ITS (Shannon Diversity)"      -> col 3  Farm, Site E:
#                                             5.5 ± 0.8 -> 5.6 ± 0.8
#                                          -> col 4  Natural, Site F (the
#                                             first "6.0 ± 0.3" in the row):
#                                             6.0 ± 0.3 -> 6.0 ± 0.2
#                                          -> col 9  Farm, Site G:
#                                             5.0 ± 0.5 -> 5.1 ± 0.5

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $oldText, $newText) {
    # Cell.Range's own Find appears to operate against the whole story
    # rather than being clipped to the cell, so re-wrap the cell's
    # [Start,End) offsets in a fresh Document.Range before searching --
    # that keeps the replacement confined to this single cell even when
    # neighbouring cells contain the same text.
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    $scoped = $d.Range($cellRange.Start, $cellRange.End)

    $found = $scoped.Find.Execute($oldText, $true, $false, $false, $false,
                                   $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Could not find '$oldText' in cell ($row, $col)"
    }
}

Set-CellText $t 10 13 "0.016" "0.023"
Set-CellText $t 11 3  "5.5 ± 0.8" "5.6 ± 0.8"
Set-CellText $t 11 4  "6.0 ± 0.3" "6.0 ± 0.2"
Set-CellText $t 11 9  "5.0 ± 0.5" "5.1 ± 0.5"
